# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

# row -> new F-column value
$updates = @{
    3  = 614
    4  = 2176
    5  = 67
    6  = 12736
    10 = 466
    13 = 13680
    14 = 14074
    16 = 169
    18 = 42
    19 = 21
    26 = 928
    27 = 5191
    29 = 264
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
